$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the test data record values (row 8) in the shared strings table
$ws.Range("A8").Value = "ATestCAutomationC1"
$ws.Range("B8").Value = "ATestCAutomationC1"
$ws.Range("C8").Value = "Facility_E1091a1222"
$ws.Range("D8").Value = "Facility_E1091a1222"
$ws.Range("E8").Value = "Pharmacy_E1091a1222"
$ws.Range("F8").Value = "Pharmacy_E1091a1222"
$ws.Range("H8").Value = "Alignment Project E1091a1222"

# Update the active cell selection to I13
$ws.Range("I13").Select()

$wb.Save()
